# The sheet had an extraneous "Unnamed: 0" index column in column B
# (leftover from a pandas export). Remove it so every later column
# (lastname, firstname, middlename, email, cash, password, sent)
# shifts one slot to the left, and reset the "sent" flag (now column H)
# back to 0 for the two existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Unnamed: 0") - shifts C:I left to B:H.
$ws.Columns("B").Delete()

# The former "sent" column (I) is now H; reset its values to 0.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
